{"js": "// Update the division-problem worksheet: replace each \"dividend\u00f7divisor=\"\n// expression in the table with its new value, per the commit's regenerated\n// numbers. Mapping is old-text -> new-text; old values are unique in the\n// document so an exact-match search is unambiguous.\nconst replacements = [\n  [\"242\u00f72=\", \"856\u00f78=\"],\n  [\"442\u00f79=\", \"329\u00f79=\"],\n  [\"820\u00f78=\", \"834\u00f78=\"],\n  [\"955\u00f75=\", \"880\u00f74=\"],\n  [\"494\u00f78=\", \"933\u00f77=\"],\n  [\"230\u00f74=\", \"329\u00f77=\"],\n  [\"887\u00f75=\", \"602\u00f73=\"],\n  [\"738\u00f74=\", \"986\u00f73=\"],\n  [\"497\u00f77=\", \"646\u00f78=\"],\n  [\"825\u00f78=\", \"611\u00f78=\"],\n  [\"937\u00f72=\", \"491\u00f77=\"],\n  [\"606\u00f78=\", \"285\u00f73=\"],\n  [\"395\u00f75=\", \"177\u00f73=\"],\n  [\"136\u00f75=\", \"558\u00f74=\"],\n  [\"945\u00f76=\", \"843\u00f79=\"],\n  [\"717\u00f75=\", \"130\u00f72=\"],\n  [\"790\u00f72=\", \"521\u00f74=\"],\n  [\"114\u00f78=\", \"379\u00f73=\"],\n  [\"249\u00f72=\", \"418\u00f74=\"],\n  [\"370\u00f72=\", \"383\u00f72=\"],\n  [\"290\u00f75=\", \"899\u00f77=\"],\n  [\"621\u00f75=\", \"792\u00f78=\"],\n  [\"544\u00f76=\", \"696\u00f77=\"],\n  [\"803\u00f72=\", \"802\u00f77=\"],\n  [\"114\u00f72=\", \"958\u00f76=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, {\n    matchCase: true,\n    matchWholeWord: false,\n  });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the division-problem worksheet: replace each \"dividend\u00f7divisor=\"\n# expression in the table with its new value, per the commit's regenerated\n# numbers. Mapping is old-text -> new-text; old values are unique in the\n# document so an exact-match Find/Replace is unambiguous.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"242\u00f72=\", \"856\u00f78=\"),\n    @(\"442\u00f79=\", \"329\u00f79=\"),\n    @(\"820\u00f78=\", \"834\u00f78=\"),\n    @(\"955\u00f75=\", \"880\u00f74=\"),\n    @(\"494\u00f78=\", \"933\u00f77=\"),\n    @(\"230\u00f74=\", \"329\u00f77=\"),\n    @(\"887\u00f75=\", \"602\u00f73=\"),\n    @(\"738\u00f74=\", \"986\u00f73=\"),\n    @(\"497\u00f77=\", \"646\u00f78=\"),\n    @(\"825\u00f78=\", \"611\u00f78=\"),\n    @(\"937\u00f72=\", \"491\u00f77=\"),\n    @(\"606\u00f78=\", \"285\u00f73=\"),\n    @(\"395\u00f75=\", \"177\u00f73=\"),\n    @(\"136\u00f75=\", \"558\u00f74=\"),\n    @(\"945\u00f76=\", \"843\u00f79=\"),\n    @(\"717\u00f75=\", \"130\u00f72=\"),\n    @(\"790\u00f72=\", \"521\u00f74=\"),\n    @(\"114\u00f78=\", \"379\u00f73=\"),\n    @(\"249\u00f72=\", \"418\u00f74=\"),\n    @(\"370\u00f72=\", \"383\u00f72=\"),\n    @(\"290\u00f75=\", \"899\u00f77=\"),\n    @(\"621\u00f75=\", \"792\u00f78=\"),\n    @(\"544\u00f76=\", \"696\u00f77=\"),\n    @(\"803\u00f72=\", \"802\u00f77=\"),\n    @(\"114\u00f72=\", \"958\u00f76=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
